$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.337.72"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.621.76"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.249"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "1.847.46"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "1.621.79"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "26.345.05"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0515"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.57%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("D36").Value = "1.160.60"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.495"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.783"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "1.759.16"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.409"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "
